# Add a new user row to the "Users" sheet and a new "Tables" sheet
# listing the timetable slots for that user, as per the upstream commit.
#
# NOTE: this host's `.Value` getter on a freshly-touched Range echoes the
# property's reflection signature instead of the real value, so every
# read/write below goes through `.Value2` instead (write behaviour is the
# same for plain text/number content).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Users sheet: append row 2 with the new user's data.
# ---------------------------------------------------------------------
$users = $wb.Worksheets.Item("Users")

# id-like numeric strings need to be forced to text so they round-trip
# as shared strings instead of numbers (matches the id/actionId columns
# which are stored as text elsewhere in the workbook).
$users.Range("A2").NumberFormat = "@"
$users.Range("A2").Value2 = "674723828"

$users.Range("B2").Value2 = "Maria"
$users.Range("C2").Value2 = "Belyaeva"
$users.Range("D2").Value2 = "belyaevaMar"
$users.Range("E2").Value2 = "Мария"

$users.Range("F2").NumberFormat = "@"
$users.Range("F2").Value2 = "0"

[void]$users.Range("A2:B3").Select()

# ---------------------------------------------------------------------
# 2) New "Tables" sheet (placed after "Users"), listing the user's
#    scheduled talks.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$tables = $wb.Worksheets.Add($null, $lastSheet)
$tables.Name = "Tables"

$tables.Range("A1").NumberFormat = "@"
$tables.Range("A1").Value2 = "userId"
$tables.Range("B1").NumberFormat = "@"
$tables.Range("B1").Value2 = "table"

# Reuse the already-text "id" value from the Users sheet so it is stored
# as a plain string (no extra number formatting) in the new sheet too.
$userIdText = $users.Range("A2").Value2

$tables.Range("A2").Value2 = $userIdText
$tables.Range("B2").Value2 = "16:45-17:30 — СДИ Софт (Информационная безопасность)"

$tables.Range("A3").Value2 = $userIdText
$tables.Range("B3").Value2 = "13:00-15:00 — Цифровая Россия (Зал Edison)"

# ---------------------------------------------------------------------
# 3) Update the remembered selection on the other sheets (A2:F2 -> A2:B3)
#    and make "Tables" the active sheet/tab with a 120% zoom, matching
#    the upstream diff.
# ---------------------------------------------------------------------
$admin = $wb.Worksheets.Item("Admin")
[void]$admin.Range("A2:B3").Select()

$speakers = $wb.Worksheets.Item("Speakers")
[void]$speakers.Range("A2:B3").Select()

$questions = $wb.Worksheets.Item("Questions")
[void]$questions.Range("A2:B3").Select()

$tables.Activate()
[void]$tables.Range("A2:B3").Select()
$excel.ActiveWindow.Zoom = 120
